# Auto-generated edit script: updates currentAveragePrice / LeveProfit
# columns (H-N) across several Leve-profit worksheets, per scheduled
# market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 821.0909
$ws.Range("I32").Value = 793
$ws.Range("K32").Value = 793
$ws.Range("M32").Value = -467
$ws.Range("H99").Value = 967.36365
$ws.Range("I99").Value = 880.5
$ws.Range("J99").Value = 1199
$ws.Range("K99").Value = 2641.5
$ws.Range("L99").Value = 3597
$ws.Range("M99").Value = -1143.5
$ws.Range("N99").Value = -6593
$ws.Range("H101").Value = 400.7
$ws.Range("I101").Value = 390.5
$ws.Range("J101").Value = 416
$ws.Range("K101").Value = 1171.5
$ws.Range("L101").Value = 1248
$ws.Range("M101").Value = 450.5
$ws.Range("N101").Value = -4492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4111.1055
$ws.Range("I32").Value = 1606.4
$ws.Range("K32").Value = 1606.4
$ws.Range("M32").Value = -1319.4
$ws.Range("H37").Value = 11666.667
$ws.Range("H74").Value = 4499.3335
$ws.Range("I74").Value = 4499
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 4499
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -3625
$ws.Range("N74").Value = -6248
$ws.Range("H77").Value = 4499.3335
$ws.Range("I77").Value = 4499
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 22495
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -18127
$ws.Range("N77").Value = -31236
$ws.Range("H122").Value = 3249.5
$ws.Range("I122").Value = 2999.3333
$ws.Range("K122").Value = 8997.999899999999
$ws.Range("M122").Value = -6547.999899999999
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H134").Value = 57499.332
$ws.Range("J134").Value = 57499.332
$ws.Range("L134").Value = 57499.332
$ws.Range("N134").Value = -67639.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 2324.182
$ws.Range("I31").Value = 1889.8
$ws.Range("K31").Value = 1889.8
$ws.Range("M31").Value = -1594.8
$ws.Range("H34").Value = 2324.182
$ws.Range("I34").Value = 1889.8
$ws.Range("K34").Value = 1889.8
$ws.Range("M34").Value = -1687.8
$ws.Range("H132").Value = 5366
$ws.Range("I132").Value = 5398.6665
$ws.Range("K132").Value = 16195.9995
$ws.Range("M132").Value = -13665.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 71.666664
$ws.Range("I12").Value = 25.833334
$ws.Range("J12").Value = 117.5
$ws.Range("K12").Value = 77.50000199999999
$ws.Range("L12").Value = 352.5
$ws.Range("M12").Value = 95.49999800000001
$ws.Range("N12").Value = -698.5
$ws.Range("I107").Value = 990
$ws.Range("J107").Value = 1332.6666
$ws.Range("K107").Value = 2970
$ws.Range("L107").Value = 3997.9998
$ws.Range("M107").Value = -1050
$ws.Range("N107").Value = -7837.9998
$ws.Range("H140").Value = 2620.7144
$ws.Range("I140").Value = 2058
$ws.Range("K140").Value = 6174
$ws.Range("M140").Value = -994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 30831.555
$ws.Range("I7").Value = 34823.5
$ws.Range("J7").Value = 27638
$ws.Range("K7").Value = 34823.5
$ws.Range("L7").Value = 27638
$ws.Range("M7").Value = -34711.5
$ws.Range("N7").Value = -27862
$ws.Range("H40").Value = 4650
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -1664
$ws.Range("N40").Value = -7772
$ws.Range("H61").Value = 3999
$ws.Range("I61").Value = 3999
$ws.Range("K61").Value = 3999
$ws.Range("M61").Value = -3797
$ws.Range("H113").Value = 3999
$ws.Range("I113").Value = 3999
$ws.Range("K113").Value = 3999
$ws.Range("M113").Value = -1829
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 30831.555
$ws.Range("I126").Value = 34823.5
$ws.Range("J126").Value = 27638
$ws.Range("K126").Value = 104470.5
$ws.Range("L126").Value = 82914
$ws.Range("M126").Value = -102000.5
$ws.Range("N126").Value = -87854
$ws.Range("H132").Value = 3399.8
$ws.Range("I132").Value = 3399.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10199.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7669.400000000001
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 3692.7693
$ws.Range("I136").Value = 2951.75
$ws.Range("J136").Value = 4878.4
$ws.Range("K136").Value = 8855.25
$ws.Range("L136").Value = 14635.2
$ws.Range("M136").Value = -6305.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1000
$ws.Range("I81").Value = 1000
$ws.Range("K81").Value = 2000
$ws.Range("M81").Value = -939
$ws.Range("H84").Value = 1000
$ws.Range("I84").Value = 1000
$ws.Range("K84").Value = 10000
$ws.Range("M84").Value = -4696
$ws.Range("H107").Value = 1275.2307
$ws.Range("I107").Value = 908.6667
$ws.Range("J107").Value = 2100
$ws.Range("K107").Value = 2726.0001
$ws.Range("L107").Value = 6300
$ws.Range("M107").Value = -806.0001000000002
$ws.Range("N107").Value = -10140
$ws.Range("H113").Value = 290.8
$ws.Range("I113").Value = 326
$ws.Range("K113").Value = 978
$ws.Range("M113").Value = 1192
$ws.Range("H126").Value = 33161.19
$ws.Range("I126").Value = 34211.812
$ws.Range("J126").Value = 29799.2
$ws.Range("K126").Value = 102635.436
$ws.Range("L126").Value = 89397.60000000001
$ws.Range("M126").Value = -100165.436
$ws.Range("N126").Value = -94337.60000000001
$ws.Range("H136").Value = 9239
$ws.Range("I136").Value = 8849.25
$ws.Range("K136").Value = 26547.75
$ws.Range("M136").Value = -23997.75
